$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("B2").Value = 8.50525
$ws.Range("C2").Value = 17.20945
$ws.Range("D2").Value = 22.17555
$ws.Range("E2").Value = 0.12255
$ws.Range("F2").Value = 2935.1395
$ws.Range("G2").Value = 2689.01645
$ws.Range("H2").Value = 246.12305
$ws.Range("I2").Value = 246.12305
$ws.Range("K2").Value = 5325.393399999999
$ws.Range("L2").Value = 2689.035
$ws.Range("M2").Value = 2636.3584
$ws.Range("N2").Value = 395.8867999999999
$ws.Range("O2").Value = 2240.4718
$ws.Range("B3").Value = 9.336
$ws.Range("C3").Value = 24.568
$ws.Range("F3").Value = 3244.394
$ws.Range("G3").Value = 2902.063
$ws.Range("H3").Value = 342.332
$ws.Range("I3").Value = 211.476
$ws.Range("J3").Value = 130.856
$ws.Range("K3").Value = 3856.4258
$ws.Range("L3").Value = 2902.004
$ws.Range("M3").Value = 954.4217999999998
$ws.Range("N3").Value = 369.761
$ws.Range("O3").Value = 584.6609999999999
$ws.Range("B4").Value = 11.824
$ws.Range("C4").Value = 30.715
$ws.Range("F4").Value = 3597.835
$ws.Range("G4").Value = 3427.495
$ws.Range("H4").Value = 170.34
$ws.Range("I4").Value = 170.34
$ws.Range("K4").Value = 3570.5422
$ws.Range("L4").Value = 3427.495
$ws.Range("M4").Value = 143.0472
$ws.Range("N4").Value = 143.0472

$ws = $wb.Worksheets.Item(2)
$ws.Range("B2").Value = 8.720699999999999
$ws.Range("C2").Value = 18.2932
$ws.Range("D2").Value = 20.2373
$ws.Range("E2").Value = 0.15595
$ws.Range("F2").Value = 2999.0962
$ws.Range("G2").Value = 2741.5145
$ws.Range("H2").Value = 257.5818
$ws.Range("I2").Value = 257.5818
$ws.Range("K2").Value = 3299.7192
$ws.Range("L2").Value = 2741.534
$ws.Range("M2").Value = 558.1852
$ws.Range("N2").Value = 249.3912
$ws.Range("O2").Value = 308.794
$ws.Range("B3").Value = 10.121
$ws.Range("C3").Value = 26.355
$ws.Range("F3").Value = 3325.352
$ws.Range("G3").Value = 3061.865
$ws.Range("H3").Value = 263.487
$ws.Range("I3").Value = 180.566
$ws.Range("J3").Value = 82.922
$ws.Range("K3").Value = 3221.406
$ws.Range("L3").Value = 3061.917
$ws.Range("M3").Value = 159.489
$ws.Range("N3").Value = 159.489
$ws.Range("B4").Value = 12.378
$ws.Range("C4").Value = 28.967
$ws.Range("F4").Value = 3618.86
$ws.Range("G4").Value = 3424.175
$ws.Range("H4").Value = 194.685
$ws.Range("I4").Value = 194.685
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 3490.7982
$ws.Range("L4").Value = 3424.175
$ws.Range("M4").Value = 66.62320000000001
$ws.Range("N4").Value = 66.62320000000001

$ws = $wb.Worksheets.Item(3)
$ws.Range("B2").Value = 8.89795
$ws.Range("C2").Value = 17.95005
$ws.Range("D2").Value = 18.71045
$ws.Range("E2").Value = 0.1394
$ws.Range("F2").Value = 2977.2507
$ws.Range("G2").Value = 2737.2133
$ws.Range("H2").Value = 240.03745
$ws.Range("I2").Value = 240.03745
$ws.Range("K2").Value = 3574.756999999999
$ws.Range("L2").Value = 2737.232
$ws.Range("M2").Value = 837.525
$ws.Range("N2").Value = 342.9553999999999
$ws.Range("O2").Value = 494.5696
$ws.Range("B3").Value = 10.126
$ws.Range("C3").Value = 26.322
$ws.Range("F3").Value = 3304.498
$ws.Range("G3").Value = 3061.199
$ws.Range("H3").Value = 243.299
$ws.Range("I3").Value = 159.898
$ws.Range("J3").Value = 83.401
$ws.Range("K3").Value = 3304.4606
$ws.Range("L3").Value = 3061.222
$ws.Range("M3").Value = 243.2386
$ws.Range("N3").Value = 243.2386
$ws.Range("B4").Value = 12.378
$ws.Range("C4").Value = 28.967
$ws.Range("F4").Value = 3618.86
$ws.Range("G4").Value = 3424.175
$ws.Range("H4").Value = 194.685
$ws.Range("I4").Value = 194.685
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 3543.8584
$ws.Range("L4").Value = 3424.175
$ws.Range("M4").Value = 119.6834
$ws.Range("N4").Value = 119.6834

$ws = $wb.Worksheets.Item(4)
$ws.Range("B2").Value = 9.119949999999999
$ws.Range("C2").Value = 18.2508
$ws.Range("D2").Value = 22.10749999999999
$ws.Range("E2").Value = 0.15595
$ws.Range("F2").Value = 3058.390549999999
$ws.Range("G2").Value = 2800.1548
$ws.Range("H2").Value = 258.23585
$ws.Range("I2").Value = 258.23585
$ws.Range("K2").Value = 2909.3964
$ws.Range("L2").Value = 2800.176
$ws.Range("M2").Value = 109.2204
$ws.Range("N2").Value = 109.2204
$ws.Range("B3").Value = 10.121
$ws.Range("C3").Value = 26.355
$ws.Range("F3").Value = 3345.42
$ws.Range("G3").Value = 3061.865
$ws.Range("H3").Value = 283.555
$ws.Range("I3").Value = 200.633
$ws.Range("J3").Value = 82.922
$ws.Range("K3").Value = 3141.1456
$ws.Range("L3").Value = 3061.917
$ws.Range("M3").Value = 79.2286
$ws.Range("N3").Value = 79.2286
$ws.Range("B4").Value = 12.378
$ws.Range("C4").Value = 28.967
$ws.Range("F4").Value = 3618.86
$ws.Range("G4").Value = 3424.175
$ws.Range("H4").Value = 194.685
$ws.Range("I4").Value = 194.685
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 3447.0206
$ws.Range("L4").Value = 3424.175
$ws.Range("M4").Value = 22.8456
$ws.Range("N4").Value = 22.8456

$ws = $wb.Worksheets.Item(5)
$ws.Range("B2").Value = 8.877549999999999
$ws.Range("C2").Value = 17.1503
$ws.Range("D2").Value = 21.6828
$ws.Range("E2").Value = 0.04995000000000001
$ws.Range("F2").Value = 2990.318849999999
$ws.Range("G2").Value = 2725.606150000001
$ws.Range("H2").Value = 264.71265
$ws.Range("I2").Value = 264.71265
$ws.Range("K2").Value = 4044.4486
$ws.Range("L2").Value = 2725.629
$ws.Range("M2").Value = 1318.8196
$ws.Range("N2").Value = 169.79
$ws.Range("O2").Value = 1149.0296
$ws.Range("B3").Value = 10.275
$ws.Range("C3").Value = 23.599
$ws.Range("F3").Value = 3238.849
$ws.Range("G3").Value = 2973.33
$ws.Range("H3").Value = 265.519
$ws.Range("I3").Value = 229.315
$ws.Range("J3").Value = 36.204
$ws.Range("K3").Value = 3566.970400000001
$ws.Range("L3").Value = 2973.351
$ws.Range("M3").Value = 593.6193999999999
$ws.Range("N3").Value = 141.2584
$ws.Range("O3").Value = 452.361
$ws.Range("B4").Value = 11.594
$ws.Range("C4").Value = 23.812
$ws.Range("D4").Value = 8.042999999999999
$ws.Range("F4").Value = 3458.482
$ws.Range("G4").Value = 3194.709
$ws.Range("H4").Value = 263.773
$ws.Range("I4").Value = 187.134
$ws.Range("J4").Value = 76.639
$ws.Range("K4").Value = 3526.0414
$ws.Range("L4").Value = 3194.709
$ws.Range("M4").Value = 331.3324
$ws.Range("N4").Value = 80.1902
$ws.Range("O4").Value = 251.1422
